# Generate Report for handback
# Update the "Correspond Handoff Datetime" (column D) and
# "Correspond Handback DateTime" (column G) values for the rows that were
# just handed back, on both the "zh-cn" and "de-de" report sheets.

$wb = $excel.ActiveWorkbook

# --- zh-cn sheet: row 3 (acb8a402...zh-cn.xlf) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D3").Value = "2016-02-15 03:46:15"
$wsZhCn.Range("G3").Value = "2016-02-15 03:47:18"

# --- de-de sheet: row 3 (acb8a402...de-de.xlf) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D3").Value = "2016-02-15 03:46:32"
$wsDeDe.Range("G3").Value = "2016-02-15 03:47:42"
